$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the RANDBETWEEN formulas in B2:B22 with the static value 11
$ws.Range("B2:B22").Value = 11

# Update the selection to match the new range
$ws.Range("B2:B22").Select()
